# The workbook gained one new data row. A new record was inserted at
# worksheet row 188 ("Feria Lagunitas de Puerto Montt" / Zanahoria data),
# which pushed all the previously-existing rows 188-270 down by one
# position (to 189-271), growing the used range from A1:R270 to A1:R271.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 188, shifting rows 188:270 down to 189:271.
$ws.Rows("188:188").Insert()

# Populate the newly inserted row 188 with the new record's data.
$ws.Cells.Item(188, 1).Value  = 4
$ws.Cells.Item(188, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(188, 3).Value  = "Los Lagos"
$ws.Cells.Item(188, 4).Value  = 44553
$ws.Cells.Item(188, 5).Value  = 10
$ws.Cells.Item(188, 6).Value  = 100114013
$ws.Cells.Item(188, 7).Value  = "Zanahoria"
$ws.Cells.Item(188, 8).Value  = "Sin especificar"
$ws.Cells.Item(188, 9).Value  = "Primera"
$ws.Cells.Item(188, 10).Value = 300
$ws.Cells.Item(188, 11).Value = 11500
$ws.Cells.Item(188, 12).Value = 12000
$ws.Cells.Item(188, 13).Value = 11750
$ws.Cells.Item(188, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(188, 15).Value = "Región de Ñuble"
$ws.Cells.Item(188, 16).Value = 588
$ws.Cells.Item(188, 17).Value = 20
$ws.Cells.Item(188, 18).Value = "Hortaliza"

# Make sure the D column (date) of the new row uses the same date number
# format style as the rest of the date column.
$ws.Range("D188").NumberFormat = $ws.Range("D189").NumberFormat
